# Auto-generated edit script: applies the "Updated symbol list on Mon Jan 23 06:56:52 UTC 2023 with GitHub Actions" diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: Coin/Link columns (B, C) are plain text and can be set directly.
# Price/Volume columns (D, E) contain numeric- and percent-looking text that Excel
# would otherwise auto-convert to Number/Percentage; force them to stay as literal
# text (matching the source inline-string cells) by applying a Text number format
# before the write, then restoring the "Normal" style so no stray formatting sticks.
$updates = @(
    @{ Cell = "D2"; Value = "305.21" }
    @{ Cell = "E2"; Value = "1.16%" }
    @{ Cell = "D3"; Value = "36.05" }
    @{ Cell = "E3"; Value = "-3.84%" }
    @{ Cell = "D4"; Value = "5.124" }
    @{ Cell = "E4"; Value = "2.31%" }
    @{ Cell = "D5"; Value = "0.07875" }
    @{ Cell = "E5"; Value = "0.29%" }
    @{ Cell = "D6"; Value = "2.183" }
    @{ Cell = "E6"; Value = "-3.62%" }
    @{ Cell = "D7"; Value = "7.935" }
    @{ Cell = "E7"; Value = "-1.15%" }
    @{ Cell = "B8"; Value = "MXToken" }
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D8"; Value = "0.9194" }
    @{ Cell = "E8"; Value = "1.05%" }
    @{ Cell = "B9"; Value = "LiechtensteinCryptoassetsExchange" }
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" }
    @{ Cell = "D9"; Value = "0.09689" }
    @{ Cell = "E9"; Value = "4.66%" }
    @{ Cell = "B10"; Value = "WazirX" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "D10"; Value = "0.1875" }
    @{ Cell = "E10"; Value = "-0.61%" }
    @{ Cell = "B11"; Value = "MandalaExchangeToken" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D11"; Value = "0.08688" }
    @{ Cell = "E11"; Value = "2.45%" }
    @{ Cell = "B12"; Value = "BitrueCoin" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D12"; Value = "0.03482" }
    @{ Cell = "E12"; Value = "-1.50%" }
    @{ Cell = "B13"; Value = "BitMartToken" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D13"; Value = "0.09931" }
    @{ Cell = "E13"; Value = "-0.16%" }
    @{ Cell = "B14"; Value = "BitForexToken" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D14"; Value = "0.001448" }
    @{ Cell = "E14"; Value = "-2.58%" }
    @{ Cell = "B15"; Value = "TigerCash" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D15"; Value = "0.005700" }
    @{ Cell = "E15"; Value = "1.25%" }
    @{ Cell = "B16"; Value = "LEO" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D16"; Value = "3.460" }
    @{ Cell = "E16"; Value = "-0.33%" }
    @{ Cell = "B17"; Value = "GateToken" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Cell = "D17"; Value = "4.098" }
    @{ Cell = "E17"; Value = "2.01%" }
    @{ Cell = "D18"; Value = "2.465" }
    @{ Cell = "E18"; Value = "14.49%" }
    @{ Cell = "E19"; Value = "-1.09%" }
    @{ Cell = "D20"; Value = "0.1299" }
    @{ Cell = "E20"; Value = "-0.63%" }
    @{ Cell = "D21"; Value = "4.844" }
    @{ Cell = "E21"; Value = "1.56%" }
    @{ Cell = "E22"; Value = "-0.02%" }
    @{ Cell = "D23"; Value = "0.04543" }
    @{ Cell = "E23"; Value = "-2.29%" }
    @{ Cell = "E24"; Value = "14.40%" }
    @{ Cell = "D25"; Value = "0.001231" }
    @{ Cell = "E25"; Value = "0.26%" }
    @{ Cell = "D26"; Value = "0.0001400" }
    @{ Cell = "E26"; Value = "7.79%" }
    @{ Cell = "D27"; Value = "0.0004747" }
    @{ Cell = "E27"; Value = "0.07%" }
    @{ Cell = "E39"; Value = "4.61%" }
    @{ Cell = "D40"; Value = "0.04785" }
    @{ Cell = "E40"; Value = "0.96%" }
    @{ Cell = "D41"; Value = "0.007704" }
    @{ Cell = "E41"; Value = "-1.95%" }
    @{ Cell = "D42"; Value = "0.1399" }
    @{ Cell = "E42"; Value = "0.53%" }
    @{ Cell = "D43"; Value = "0.007732" }
    @{ Cell = "E43"; Value = "1.04%" }
    @{ Cell = "D44"; Value = "0.002230" }
    @{ Cell = "E44"; Value = "0.06%" }
    @{ Cell = "D45"; Value = "0.01100" }
    @{ Cell = "E45"; Value = "7.57%" }
    @{ Cell = "D46"; Value = "0.00006393" }
    @{ Cell = "E46"; Value = "5.59%" }
    @{ Cell = "E47"; Value = "0.09%" }
    @{ Cell = "D48"; Value = "0.0005797" }
    @{ Cell = "E48"; Value = "-0.06%" }
    @{ Cell = "D49"; Value = "24.46" }
    @{ Cell = "E49"; Value = "182.08%" }
    @{ Cell = "D50"; Value = "0.001999" }
    @{ Cell = "E50"; Value = "-25.59%" }
    @{ Cell = "D51"; Value = "0.00002099" }
    @{ Cell = "E51"; Value = "0.09%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

